# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly (re)computed "K" values (column G) for each data row, keyed by row number.
$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 1
    24 = 0
    25 = 0
    26 = 1
    27 = 0
    28 = 0
    29 = 1
    30 = 1
    31 = 0
    32 = 0
    33 = 2
    34 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
